# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets to reflect the latest generated data.
#   F3: 47 -> 49
#   F4: 0  -> 4

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 49
    $ws.Range("F4").Value = 4
}
